# Delete All Unnecessary excel files
# - Update row 2's file_date (C2) from "01-01" to "02-03" (formula in B2 recalculates)
# - Delete row 3 entirely (جيانا/رند sample row for "رند")
# - Update selection to C3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the file_date value in C2; formula in B2 will recalc automatically
$ws.Range("C2").Value = "02-03"

# Delete entire row 3 (shifts nothing below it up, but removes the row)
$ws.Rows("3:3").Delete()

# Update the active selection as per the diff
$ws.Range("C3").Select()

$wb.Save()
